# Update the dated header line.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-03-22 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-03-23 Sunday", 2)

# Update the division-problem answers in the single table on the page.
# Addressed by (row, column) rather than text search because one of the
# old values ("31÷3=10, 1") appears twice in the table.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "55÷8=6, 7"
$t.Cell(1, 2).Range.Text = "73÷3=24, 1"
$t.Cell(1, 3).Range.Text = "48÷5=9, 3"
$t.Cell(1, 4).Range.Text = "66÷8=8, 2"
$t.Cell(1, 5).Range.Text = "95÷9=10, 5"

$t.Cell(5, 1).Range.Text = "15÷8=1, 7"
$t.Cell(5, 2).Range.Text = "78÷3=26, 0"
$t.Cell(5, 3).Range.Text = "16÷5=3, 1"
$t.Cell(5, 4).Range.Text = "63÷4=15, 3"
$t.Cell(5, 5).Range.Text = "64÷4=16, 0"

$t.Cell(9, 1).Range.Text = "67÷5=13, 2"
$t.Cell(9, 2).Range.Text = "21÷6=3, 3"
$t.Cell(9, 3).Range.Text = "90÷7=12, 6"
$t.Cell(9, 4).Range.Text = "45÷3=15, 0"
$t.Cell(9, 5).Range.Text = "62÷5=12, 2"

$t.Cell(13, 1).Range.Text = "53÷3=17, 2"
$t.Cell(13, 2).Range.Text = "66÷6=11, 0"
$t.Cell(13, 3).Range.Text = "57÷8=7, 1"
$t.Cell(13, 4).Range.Text = "67÷7=9, 4"
$t.Cell(13, 5).Range.Text = "46÷4=11, 2"

$t.Cell(17, 1).Range.Text = "93÷4=23, 1"
$t.Cell(17, 2).Range.Text = "35÷4=8, 3"
$t.Cell(17, 3).Range.Text = "24÷8=3, 0"
$t.Cell(17, 4).Range.Text = "90÷5=18, 0"
$t.Cell(17, 5).Range.Text = "15÷9=1, 6"
